$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new numeric-looking text (e.g. "1.00") must stay stored as TEXT
# rather than being auto-coerced to a number by Excel's input parser. We
# temporarily mark each one Text-formatted before assigning the literal string.
$textSafeAddrs = @(
    "D5", "D6", "D10", "D12", "D15", "D16", "D19", "D20",
    "D21", "D22", "D23", "D24", "D25", "D29", "D31", "D32",
    "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41",
    "D43", "D44", "D45", "D46", "D49", "D50", "D51"
)
foreach ($addr in $textSafeAddrs) {
    $ws.Range($addr).NumberFormat = "@"
}

# --- Row-by-row value updates (Coin, Link, Price, Volume(1h)) ---
$ws.Range("D2").Value = '63.343.37'
$ws.Range("E2").Value = '  +5.74%  '

$ws.Range("D3").Value = '3.411.06'
$ws.Range("E3").Value = '  +6.84%  '

$ws.Range("E4").Value = '  +0.04%  '

$ws.Range("D5").Value = '576.95'
$ws.Range("E5").Value = '  +7.48%  '

$ws.Range("D6").Value = '155.55'
$ws.Range("E6").Value = '  +7.17%  '

$ws.Range("E7").Value = '  -0.01%  '

$ws.Range("D8").Value = '3.419.61'
$ws.Range("E8").Value = '  +6.91%  '

$ws.Range("E9").Value = '  +0.37%  '

$ws.Range("D10").Value = '7.54'
$ws.Range("E10").Value = '  +3.11%  '

$ws.Range("E11").Value = '  +7.86%  '

$ws.Range("D12").Value = '0.436'
$ws.Range("E12").Value = '  +0.44%  '

$ws.Range("D13").Value = '3.996.46'
$ws.Range("E13").Value = '  +6.86%  '

$ws.Range("E14").Value = '  -0.37%  '

$ws.Range("D15").Value = '0.0000185'
$ws.Range("E15").Value = '  +7.74%  '

$ws.Range("D16").Value = '27.09'
$ws.Range("E16").Value = '  +5.44%  '

$ws.Range("D17").Value = '63.544.65'
$ws.Range("E17").Value = '  +6.03%  '

$ws.Range("D18").Value = '3.355.83'
$ws.Range("E18").Value = '  +4.14%  '

$ws.Range("D19").Value = '6.40'
$ws.Range("E19").Value = '  +2.48%  '

$ws.Range("D20").Value = '14.16'
$ws.Range("E20").Value = '  +6.88%  '

$ws.Range("D21").Value = '8.47'
$ws.Range("E21").Value = '  +3.59%  '

$ws.Range("D22").Value = '390.63'
$ws.Range("E22").Value = '  +5.77%  '

$ws.Range("D23").Value = '1.00'
$ws.Range("E23").Value = '  -0.02%  '

$ws.Range("D24").Value = '0.537'
$ws.Range("E24").Value = '  +2.71%  '

$ws.Range("D25").Value = '71.34'
$ws.Range("E25").Value = '  +2.51%  '

$ws.Range("E26").Value = '  +21.48%  '

$ws.Range("E27").Value = '  +11.13%  '

$ws.Range("E28").Value = '  +7.08%  '

$ws.Range("D29").Value = '1.00'
$ws.Range("E29").Value = '  +0.80%  '

$ws.Range("E30").Value = '  +8.06%  '

$ws.Range("D31").Value = '6.54'
$ws.Range("E31").Value = '  +7.51%  '

$ws.Range("D32").Value = '1.34'
$ws.Range("E32").Value = '  +13.70%  '

$ws.Range("E33").Value = '  +8.41%  '

$ws.Range("D34").Value = '23.33'
$ws.Range("E34").Value = '  +3.95%  '

$ws.Range("D35").Value = '0.998'
$ws.Range("E35").Value = '  -0.08%  '

$ws.Range("D36").Value = '6.74'
$ws.Range("E36").Value = '  +2.66%  '

$ws.Range("D37").Value = '1.50'
$ws.Range("E37").Value = '  +9.93%  '

$ws.Range("D38").Value = '158.58'
$ws.Range("E38").Value = '  +0.16%  '

$ws.Range("B39").Value = 'Hedera'
$ws.Range("C39").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D39").Value = '0.0777'
$ws.Range("E39").Value = '  +9.82%  '

$ws.Range("B40").Value = 'Stacks'
$ws.Range("C40").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D40").Value = '1.89'
$ws.Range("E40").Value = '  +11.78%  '

$ws.Range("B41").Value = 'EnergySwap'
$ws.Range("C41").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D41").Value = '27.72'
$ws.Range("E41").Value = '  +5.00%  '

$ws.Range("D42").Value = '2.931.26'
$ws.Range("E42").Value = '  +5.33%  '

$ws.Range("D43").Value = '0.0320'
$ws.Range("E43").Value = '  +3.00%  '

$ws.Range("D44").Value = '0.764'
$ws.Range("E44").Value = '  +6.51%  '

$ws.Range("D45").Value = '41.57'
$ws.Range("E45").Value = '  +4.23%  '

$ws.Range("D46").Value = '4.33'
$ws.Range("E46").Value = '  +2.72%  '

$ws.Range("E47").Value = '  +10.15%  '

$ws.Range("D48").Value = '3.459.32'
$ws.Range("E48").Value = '  +6.99%  '

$ws.Range("D49").Value = '22.46'
$ws.Range("E49").Value = '  +9.00%  '

$ws.Range("B50").Value = 'Cosmos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D50").Value = '6.36'
$ws.Range("E50").Value = '  +3.41%  '

$ws.Range("B51").Value = 'Bittensor'
$ws.Range("C51").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D51").Value = '295.50'
$ws.Range("E51").Value = '  +12.73%  '

# Restore the default "Normal" style on the cells we forced to Text format above,
# so the saved workbook carries no leftover explicit number-format on them.
foreach ($addr in $textSafeAddrs) {
    $ws.Range($addr).Style = "Normal"
}
